# This script applies the "Updated symbol list" data refresh to the
# cryptos worksheet: updated Price (D) / Volume(1h) (E) figures for most
# rows, plus a swap of the BitrueCoin / MandalaExchangeToken rows
# (Coin name + Link + Price + Volume(1h) for rows 11 and 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin name / link columns) - no special numeric
# handling required since these values never look like numbers.
$textUpdates = @{
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price / Volume(1h) updates. These columns hold numbers-as-text in the
# original workbook (e.g. "5.200", "-3.12%") so the literal text/formatting
# must be preserved rather than letting Excel reinterpret them as a
# number or percentage. We force the cell to Text format before writing
# the value, then clear the explicit formatting again so the cell keeps
# the workbook's default (unstyled) look, matching the source data.
$valueUpdates = @{
    "D2"  = "321.56";     "E2"  = "-3.12%"
    "D3"  = "42.92";      "E3"  = "-5.74%"
    "D4"  = "5.200";      "E4"  = "-6.35%"
    "D5"  = "0.08175";    "E5"  = "-2.27%"
    "D6"  = "4.315";      "E6"  = "-2.60%"
    "D7"  = "1.813";      "E7"  = "-13.49%"
    "D8"  = "0.9345";     "E8"  = "-5.66%"
    "D9"  = "0.1107";     "E9"  = "-7.36%"
    "D10" = "0.1856";     "E10" = "-3.97%"
    "D11" = "0.09493";    "E11" = "-4.01%"
    "D12" = "0.04638";    "E12" = "-0.68%"
    "E13" = "-27.49%"
    "D14" = "0.1057";     "E14" = "-0.27%"
    "D15" = "0.001292";   "E15" = "0.32%"
    "D16" = "0.005858";   "E16" = "-1.11%"
    "E17" = "-1.18%"
    "E18" = "-2.18%"
    "E20" = "1.97%"
    "D21" = "0.2522";     "E21" = "-1.61%"
    "E22" = "0.72%"
    "D23" = "0.001242";   "E23" = "-3.95%"
    "D24" = "0.004293";   "E24" = "-5.46%"
    "E25" = "-7.78%"
    "D26" = "0.0002981";  "E26" = "-20.41%"
    "D38" = "0.02711";    "E38" = "0.45%"
    "D39" = "0.05555";    "E39" = "-3.37%"
    "D40" = "0.008083";   "E40" = "2.69%"
    "D41" = "0.1397";     "E41" = "-2.58%"
    "D42" = "0.006546";   "E42" = "-16.92%"
    "D43" = "0.002042";   "E43" = "0.89%"
    "D44" = "0.008254";   "E44" = "-7.65%"
    "D45" = "0.3499";     "E45" = "2.66%"
    "D46" = "0.00006926"; "E46" = "-1.84%"
    "E47" = "-0.12%"
    "D48" = "0.003338";   "E48" = "-4.37%"
    "D49" = "0.003532";   "E49" = "-0.12%"
    "D50" = "0.00002101"; "E50" = "-0.12%"
    "D51" = "0.0002001";  "E51" = "-0.12%"
}

foreach ($ref in $valueUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $valueUpdates[$ref]
    $cell.ClearFormats()
}
